# Update results for each year sheet with new server-computed values.
$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("A2").Value = 0
$ws2025.Range("B2").Value = 435.7823875228254
$ws2025.Range("E2").Value = 28829.95429482079
$ws2025.Range("G2").Value = 8095.925712661834
$ws2025.Range("I2").Value = 15999.27815958817
$ws2025.Range("L2").Value = 48700.22979492001
$ws2025.Range("M2").Value = 11286.13269883
$ws2025.Range("N2").Value = 7240.482101700985
$ws2025.Range("O2").Value = 6763.923782365135

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 412.8642372630347
$ws2030.Range("B2").Value = 4175.540191128684
$ws2030.Range("E2").Value = 44535.67533374021
$ws2030.Range("G2").Value = 8095.925712661834
$ws2030.Range("I2").Value = 36820.41621176117
$ws2030.Range("L2").Value = 48700.22979492001
$ws2030.Range("M2").Value = 17445.30966502775
$ws2030.Range("N2").Value = 9100.129533245807
$ws2030.Range("O2").Value = 7873.95248561621

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 2505.601390754976
$ws2035.Range("B2").Value = 5768.450008874052
$ws2035.Range("E2").Value = 56254.37625095449
$ws2035.Range("G2").Value = 8095.925712661834
$ws2035.Range("I2").Value = 53564.23668729185
$ws2035.Range("L2").Value = 48700.22979492001
$ws2035.Range("M2").Value = 23877.68385614976
$ws2035.Range("N2").Value = 13513.33478279869
$ws2035.Range("O2").Value = 12997.60227839026

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Value = 2505.601390754976
$ws2040.Range("B2").Value = 5768.450008874052
$ws2040.Range("E2").Value = 56254.37625095449
$ws2040.Range("G2").Value = 8095.925712661834
$ws2040.Range("I2").Value = 53564.23668729185
$ws2040.Range("L2").Value = 48700.22979492001
$ws2040.Range("M2").Value = 23877.68385614976
$ws2040.Range("N2").Value = 13513.33478279869
$ws2040.Range("O2").Value = 12997.60227839026

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 2505.601390754976
$ws2045.Range("B2").Value = 5768.450008874052
$ws2045.Range("E2").Value = 56254.37625095449
$ws2045.Range("G2").Value = 8095.925712661834
$ws2045.Range("I2").Value = 53564.23668729185
$ws2045.Range("L2").Value = 48700.22979492001
$ws2045.Range("M2").Value = 23877.68385614976
$ws2045.Range("N2").Value = 13513.33478279869
$ws2045.Range("O2").Value = 12997.60227839026

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Value = 2505.601390754976
$ws2050.Range("B2").Value = 5768.450008874052
$ws2050.Range("E2").Value = 56254.37625095449
$ws2050.Range("G2").Value = 8095.925712661834
$ws2050.Range("I2").Value = 53564.23668729185
$ws2050.Range("L2").Value = 48700.22979492001
$ws2050.Range("M2").Value = 23877.68385614976
$ws2050.Range("N2").Value = 13513.33478279869
$ws2050.Range("O2").Value = 12997.60227839026
